$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap/rotate the F:V (match detail) columns among rows that were reordered ---
# Row 12 (from original row 13)
$ws.Cells.Item(12, 6).Value = "Murcia"
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = "Cordoba"
$ws.Cells.Item(12, 9).Value = 3
$ws.Cells.Item(12, 10).Value = 2.05
$ws.Cells.Item(12, 11).Value = "01/09/2023 17:43"
$ws.Cells.Item(12, 12).Value = 2.08
$ws.Cells.Item(12, 13).Value = "02/09/2023 20:22"
$ws.Cells.Item(12, 14).Value = 3.13
$ws.Cells.Item(12, 15).Value = "01/09/2023 17:43"
$ws.Cells.Item(12, 16).Value = 3.25
$ws.Cells.Item(12, 17).Value = "02/09/2023 20:22"
$ws.Cells.Item(12, 18).Value = 3.53
$ws.Cells.Item(12, 19).Value = "01/09/2023 17:43"
$ws.Cells.Item(12, 20).Value = 3.7
$ws.Cells.Item(12, 21).Value = "02/09/2023 20:22"
$ws.Cells.Item(12, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/murcia-cordoba/8bDXxDc1/"

# Row 13 (from original row 12)
$ws.Cells.Item(13, 6).Value = "Linares"
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = "Granada CF B"
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 2.01
$ws.Cells.Item(13, 11).Value = "01/09/2023 17:43"
$ws.Cells.Item(13, 12).Value = 2.64
$ws.Cells.Item(13, 13).Value = "02/09/2023 20:21"
$ws.Cells.Item(13, 14).Value = 3.2
$ws.Cells.Item(13, 15).Value = "01/09/2023 17:43"
$ws.Cells.Item(13, 16).Value = 3.16
$ws.Cells.Item(13, 17).Value = "02/09/2023 18:33"
$ws.Cells.Item(13, 18).Value = 3.44
$ws.Cells.Item(13, 19).Value = "01/09/2023 17:43"
$ws.Cells.Item(13, 20).Value = 2.76
$ws.Cells.Item(13, 21).Value = "02/09/2023 20:21"
$ws.Cells.Item(13, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/linares-granada-cf/KYGyxXC7/"

# Row 25 (from original row 28)
$ws.Cells.Item(25, 6).Value = "Merida AD"
$ws.Cells.Item(25, 7).Value = 2
$ws.Cells.Item(25, 8).Value = "Melilla"
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 1.92
$ws.Cells.Item(25, 11).Value = "07/09/2023 09:13"
$ws.Cells.Item(25, 12).Value = 2.1
$ws.Cells.Item(25, 13).Value = "10/09/2023 08:28"
$ws.Cells.Item(25, 14).Value = 3.1
$ws.Cells.Item(25, 15).Value = "07/09/2023 09:13"
$ws.Cells.Item(25, 16).Value = 3.1
$ws.Cells.Item(25, 17).Value = "10/09/2023 10:05"
$ws.Cells.Item(25, 18).Value = 3.89
$ws.Cells.Item(25, 19).Value = "07/09/2023 09:13"
$ws.Cells.Item(25, 20).Value = 3.79
$ws.Cells.Item(25, 21).Value = "10/09/2023 08:28"
$ws.Cells.Item(25, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/merida-ad-melilla/lSX4CDlK/"

# Row 27 (from original row 25)
$ws.Cells.Item(27, 6).Value = "Atl. Madrid B"
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = "Ceuta"
$ws.Cells.Item(27, 9).Value = 3
$ws.Cells.Item(27, 10).Value = 2.24
$ws.Cells.Item(27, 11).Value = "07/09/2023 09:13"
$ws.Cells.Item(27, 12).Value = 2.18
$ws.Cells.Item(27, 13).Value = "09/09/2023 14:59"
$ws.Cells.Item(27, 14).Value = 2.94
$ws.Cells.Item(27, 15).Value = "07/09/2023 09:13"
$ws.Cells.Item(27, 16).Value = 3.1
$ws.Cells.Item(27, 17).Value = "10/09/2023 10:02"
$ws.Cells.Item(27, 18).Value = 3.18
$ws.Cells.Item(27, 19).Value = "07/09/2023 09:13"
$ws.Cells.Item(27, 20).Value = 3.57
$ws.Cells.Item(27, 21).Value = "09/09/2023 14:59"
$ws.Cells.Item(27, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/atl-madrid-ad-ceuta/QZqbgYsr/"

# Row 28 (from original row 27)
$ws.Cells.Item(28, 6).Value = "Alcoyano"
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = "UD Ibiza"
$ws.Cells.Item(28, 9).Value = 2
$ws.Cells.Item(28, 10).Value = 2.42
$ws.Cells.Item(28, 11).Value = "07/09/2023 09:13"
$ws.Cells.Item(28, 12).Value = 2.55
$ws.Cells.Item(28, 13).Value = "10/09/2023 11:42"
$ws.Cells.Item(28, 14).Value = 2.9
$ws.Cells.Item(28, 15).Value = "07/09/2023 09:13"
$ws.Cells.Item(28, 16).Value = 3
$ws.Cells.Item(28, 17).Value = "10/09/2023 10:05"
$ws.Cells.Item(28, 18).Value = 2.91
$ws.Cells.Item(28, 19).Value = "07/09/2023 09:13"
$ws.Cells.Item(28, 20).Value = 3.02
$ws.Cells.Item(28, 21).Value = "10/09/2023 11:42"
$ws.Cells.Item(28, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/alcoyano-ud-ibiza/8fcWcEZQ/"

# Row 29 (from original row 30)
$ws.Cells.Item(29, 6).Value = "San Fernando"
$ws.Cells.Item(29, 7).Value = 3
$ws.Cells.Item(29, 8).Value = "Cordoba"
$ws.Cells.Item(29, 9).Value = 1
$ws.Cells.Item(29, 10).Value = 2.63
$ws.Cells.Item(29, 11).Value = "07/09/2023 09:13"
$ws.Cells.Item(29, 12).Value = 2.83
$ws.Cells.Item(29, 13).Value = "10/09/2023 18:40"
$ws.Cells.Item(29, 14).Value = 2.95
$ws.Cells.Item(29, 15).Value = "07/09/2023 09:13"
$ws.Cells.Item(29, 16).Value = 3.06
$ws.Cells.Item(29, 17).Value = "10/09/2023 17:31"
$ws.Cells.Item(29, 18).Value = 2.62
$ws.Cells.Item(29, 19).Value = "07/09/2023 09:13"
$ws.Cells.Item(29, 20).Value = 2.65
$ws.Cells.Item(29, 21).Value = "10/09/2023 18:40"
$ws.Cells.Item(29, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/cd-san-fernando-cordoba/YJkIlAdD/"

# Row 30 (from original row 29)
$ws.Cells.Item(30, 6).Value = "Castellon"
$ws.Cells.Item(30, 7).Value = 3
$ws.Cells.Item(30, 8).Value = "CF Intercity"
$ws.Cells.Item(30, 9).Value = 1
$ws.Cells.Item(30, 10).Value = 1.77
$ws.Cells.Item(30, 11).Value = "07/09/2023 09:13"
$ws.Cells.Item(30, 12).Value = 1.56
$ws.Cells.Item(30, 13).Value = "10/09/2023 15:26"
$ws.Cells.Item(30, 14).Value = 3.24
$ws.Cells.Item(30, 15).Value = "07/09/2023 09:13"
$ws.Cells.Item(30, 16).Value = 3.8
$ws.Cells.Item(30, 17).Value = "10/09/2023 19:04"
$ws.Cells.Item(30, 18).Value = 4.39
$ws.Cells.Item(30, 19).Value = "07/09/2023 09:13"
$ws.Cells.Item(30, 20).Value = 6.5
$ws.Cells.Item(30, 21).Value = "10/09/2023 19:04"
$ws.Cells.Item(30, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/castellon-cf-intercity/U3bSbfKK/"

# Row 35 (from original row 37)
$ws.Cells.Item(35, 6).Value = "Sanluqueno"
$ws.Cells.Item(35, 7).Value = 5
$ws.Cells.Item(35, 8).Value = "Baleares"
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 2.12
$ws.Cells.Item(35, 11).Value = "14/09/2023 15:43"
$ws.Cells.Item(35, 12).Value = 1.72
$ws.Cells.Item(35, 13).Value = "17/09/2023 11:58"
$ws.Cells.Item(35, 14).Value = 3.12
$ws.Cells.Item(35, 15).Value = "14/09/2023 15:43"
$ws.Cells.Item(35, 16).Value = 3.56
$ws.Cells.Item(35, 17).Value = "17/09/2023 11:57"
$ws.Cells.Item(35, 18).Value = 3.24
$ws.Cells.Item(35, 19).Value = "14/09/2023 15:43"
$ws.Cells.Item(35, 20).Value = 5.15
$ws.Cells.Item(35, 21).Value = "17/09/2023 11:58"
$ws.Cells.Item(35, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/sanluqueno-baleares/YHKB84mn/"

# Row 37 (from original row 35)
$ws.Cells.Item(37, 6).Value = "Ceuta"
$ws.Cells.Item(37, 7).Value = 1
$ws.Cells.Item(37, 8).Value = "Alcoyano"
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 1.78
$ws.Cells.Item(37, 11).Value = "14/09/2023 15:42"
$ws.Cells.Item(37, 12).Value = 2.02
$ws.Cells.Item(37, 13).Value = "17/09/2023 11:51"
$ws.Cells.Item(37, 14).Value = 3.36
$ws.Cells.Item(37, 15).Value = "14/09/2023 15:42"
$ws.Cells.Item(37, 16).Value = 3.36
$ws.Cells.Item(37, 17).Value = "17/09/2023 11:51"
$ws.Cells.Item(37, 18).Value = 4.29
$ws.Cells.Item(37, 19).Value = "14/09/2023 15:42"
$ws.Cells.Item(37, 20).Value = 3.79
$ws.Cells.Item(37, 21).Value = "17/09/2023 11:51"
$ws.Cells.Item(37, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/ad-ceuta-alcoyano/nNFi39B6/"

# Row 38 (from original row 39)
$ws.Cells.Item(38, 6).Value = "Algeciras"
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = "Merida AD"
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 2.15
$ws.Cells.Item(38, 11).Value = "14/09/2023 15:43"
$ws.Cells.Item(38, 12).Value = 2.22
$ws.Cells.Item(38, 13).Value = "17/09/2023 14:55"
$ws.Cells.Item(38, 14).Value = 2.88
$ws.Cells.Item(38, 15).Value = "14/09/2023 15:43"
$ws.Cells.Item(38, 16).Value = 2.99
$ws.Cells.Item(38, 17).Value = "17/09/2023 17:33"
$ws.Cells.Item(38, 18).Value = 3.45
$ws.Cells.Item(38, 19).Value = "14/09/2023 15:43"
$ws.Cells.Item(38, 20).Value = 3.62
$ws.Cells.Item(38, 21).Value = "17/09/2023 15:15"
$ws.Cells.Item(38, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/algeciras-merida-ad/pbJz7iCs/"

# Row 39 (from original row 38)
$ws.Cells.Item(39, 6).Value = "Antequera"
$ws.Cells.Item(39, 7).Value = 1
$ws.Cells.Item(39, 8).Value = "Castellon"
$ws.Cells.Item(39, 9).Value = 1
$ws.Cells.Item(39, 10).Value = 3.96
$ws.Cells.Item(39, 11).Value = "14/09/2023 15:43"
$ws.Cells.Item(39, 12).Value = 3.05
$ws.Cells.Item(39, 13).Value = "17/09/2023 19:22"
$ws.Cells.Item(39, 14).Value = 3.06
$ws.Cells.Item(39, 15).Value = "14/09/2023 15:43"
$ws.Cells.Item(39, 16).Value = 3.26
$ws.Cells.Item(39, 17).Value = "17/09/2023 19:22"
$ws.Cells.Item(39, 18).Value = 1.92
$ws.Cells.Item(39, 19).Value = "14/09/2023 15:43"
$ws.Cells.Item(39, 20).Value = 2.36
$ws.Cells.Item(39, 21).Value = "17/09/2023 19:22"
$ws.Cells.Item(39, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/antequera-castellon/tQCq5Vtf/"

# Row 54 (from original row 55)
$ws.Cells.Item(54, 6).Value = "Algeciras"
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = "Sanluqueno"
$ws.Cells.Item(54, 9).Value = 1
$ws.Cells.Item(54, 10).Value = 2.29
$ws.Cells.Item(54, 11).Value = "29/09/2023 23:12"
$ws.Cells.Item(54, 12).Value = 1.9
$ws.Cells.Item(54, 13).Value = "01/10/2023 11:44"
$ws.Cells.Item(54, 14).Value = 2.89
$ws.Cells.Item(54, 15).Value = "29/09/2023 23:12"
$ws.Cells.Item(54, 16).Value = 3.22
$ws.Cells.Item(54, 17).Value = "01/10/2023 11:44"
$ws.Cells.Item(54, 18).Value = 3.14
$ws.Cells.Item(54, 19).Value = "29/09/2023 23:12"
$ws.Cells.Item(54, 20).Value = 4.54
$ws.Cells.Item(54, 21).Value = "01/10/2023 11:44"
$ws.Cells.Item(54, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/algeciras-sanluqueno/OSlaIuHA/"

# Row 55 (from original row 54)
$ws.Cells.Item(55, 6).Value = "Antequera"
$ws.Cells.Item(55, 7).Value = 3
$ws.Cells.Item(55, 8).Value = "Baleares"
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 2.01
$ws.Cells.Item(55, 11).Value = "29/09/2023 15:42"
$ws.Cells.Item(55, 12).Value = 1.8
$ws.Cells.Item(55, 13).Value = "01/10/2023 11:59"
$ws.Cells.Item(55, 14).Value = 3.04
$ws.Cells.Item(55, 15).Value = "29/09/2023 15:42"
$ws.Cells.Item(55, 16).Value = 3.49
$ws.Cells.Item(55, 17).Value = "01/10/2023 11:59"
$ws.Cells.Item(55, 18).Value = 3.64
$ws.Cells.Item(55, 19).Value = "29/09/2023 15:42"
$ws.Cells.Item(55, 20).Value = 4.62
$ws.Cells.Item(55, 21).Value = "01/10/2023 11:59"
$ws.Cells.Item(55, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/antequera-baleares/fR5AFw1T/"

# Row 84 (from original row 86)
$ws.Cells.Item(84, 6).Value = "Melilla"
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = "Antequera"
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(84, 10).Value = 2.96
$ws.Cells.Item(84, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(84, 12).Value = 3
$ws.Cells.Item(84, 13).Value = "22/10/2023 06:00"
$ws.Cells.Item(84, 14).Value = 2.91
$ws.Cells.Item(84, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(84, 16).Value = 2.91
$ws.Cells.Item(84, 17).Value = "22/10/2023 10:01"
$ws.Cells.Item(84, 18).Value = 2.44
$ws.Cells.Item(84, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(84, 20).Value = 2.59
$ws.Cells.Item(84, 21).Value = "22/10/2023 06:00"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/melilla-antequera/vBVRtyHP/"

# Row 85 (from original row 84)
$ws.Cells.Item(85, 6).Value = "CF Intercity"
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = "Baleares"
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 1.75
$ws.Cells.Item(85, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(85, 12).Value = 1.75
$ws.Cells.Item(85, 13).Value = "22/10/2023 09:56"
$ws.Cells.Item(85, 14).Value = 3.24
$ws.Cells.Item(85, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(85, 16).Value = 3.42
$ws.Cells.Item(85, 17).Value = "22/10/2023 10:01"
$ws.Cells.Item(85, 18).Value = 4.53
$ws.Cells.Item(85, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(85, 20).Value = 5.02
$ws.Cells.Item(85, 21).Value = "22/10/2023 09:56"
$ws.Cells.Item(85, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/cf-intercity-baleares/djIgxePt/"

# Row 86 (from original row 85)
$ws.Cells.Item(86, 6).Value = "Granada CF B"
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = "UD Ibiza"
$ws.Cells.Item(86, 9).Value = 1
$ws.Cells.Item(86, 10).Value = 4.18
$ws.Cells.Item(86, 11).Value = "19/10/2023 08:13"
$ws.Cells.Item(86, 12).Value = 4.16
$ws.Cells.Item(86, 13).Value = "22/10/2023 10:38"
$ws.Cells.Item(86, 14).Value = 3.2
$ws.Cells.Item(86, 15).Value = "19/10/2023 08:13"
$ws.Cells.Item(86, 16).Value = 3.14
$ws.Cells.Item(86, 17).Value = "22/10/2023 10:03"
$ws.Cells.Item(86, 18).Value = 1.85
$ws.Cells.Item(86, 19).Value = "19/10/2023 08:13"
$ws.Cells.Item(86, 20).Value = 2
$ws.Cells.Item(86, 21).Value = "22/10/2023 04:23"
$ws.Cells.Item(86, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/granada-cf-ud-ibiza/tK5Mwmql/"

# Row 96 (from original row 97)
$ws.Cells.Item(96, 6).Value = "Antequera"
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = "Malaga"
$ws.Cells.Item(96, 9).Value = 2
$ws.Cells.Item(96, 10).Value = 3.11
$ws.Cells.Item(96, 11).Value = "26/10/2023 23:12"
$ws.Cells.Item(96, 12).Value = 2.7
$ws.Cells.Item(96, 13).Value = "29/10/2023 17:46"
$ws.Cells.Item(96, 14).Value = 2.9
$ws.Cells.Item(96, 15).Value = "26/10/2023 23:12"
$ws.Cells.Item(96, 16).Value = 2.74
$ws.Cells.Item(96, 17).Value = "29/10/2023 17:46"
$ws.Cells.Item(96, 18).Value = 2.3
$ws.Cells.Item(96, 19).Value = "26/10/2023 23:12"
$ws.Cells.Item(96, 20).Value = 3.1
$ws.Cells.Item(96, 21).Value = "29/10/2023 17:46"
$ws.Cells.Item(96, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/antequera-malaga/vZDEXWwC/"

# Row 97 (from original row 96)
$ws.Cells.Item(97, 6).Value = "Cordoba"
$ws.Cells.Item(97, 7).Value = 1
$ws.Cells.Item(97, 8).Value = "Recreativo Huelva"
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 1.74
$ws.Cells.Item(97, 11).Value = "26/10/2023 23:12"
$ws.Cells.Item(97, 12).Value = 1.52
$ws.Cells.Item(97, 13).Value = "29/10/2023 17:40"
$ws.Cells.Item(97, 14).Value = 3.29
$ws.Cells.Item(97, 15).Value = "26/10/2023 23:12"
$ws.Cells.Item(97, 16).Value = 4.03
$ws.Cells.Item(97, 17).Value = "29/10/2023 17:40"
$ws.Cells.Item(97, 18).Value = 4.71
$ws.Cells.Item(97, 19).Value = "26/10/2023 23:12"
$ws.Cells.Item(97, 20).Value = 6.68
$ws.Cells.Item(97, 21).Value = "29/10/2023 17:40"
$ws.Cells.Item(97, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/cordoba-recreativo-huelva/W8A6Zh9a/"

# --- Append new rows 112-115 ---
# Row 112
$ws.Cells.Item(111, 1).Copy() | Out-Null
$ws.Cells.Item(112, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111, 5).Copy() | Out-Null
$ws.Cells.Item(112, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = "spain"
$ws.Cells.Item(112, 3).Value = "primera-rfef-group-2"
$ws.Cells.Item(112, 4).Value = "2023-2024"
$ws.Cells.Item(112, 5).Value = 45242.5
$ws.Cells.Item(112, 6).Value = "Cordoba"
$ws.Cells.Item(112, 7).Value = 3
$ws.Cells.Item(112, 8).Value = "Ceuta"
$ws.Cells.Item(112, 9).Value = 3
$ws.Cells.Item(112, 10).Value = 1.97
$ws.Cells.Item(112, 11).Value = "09/11/2023 09:13"
$ws.Cells.Item(112, 12).Value = 1.74
$ws.Cells.Item(112, 13).Value = "12/11/2023 11:58"
$ws.Cells.Item(112, 14).Value = 3.02
$ws.Cells.Item(112, 15).Value = "09/11/2023 09:13"
$ws.Cells.Item(112, 16).Value = 3.47
$ws.Cells.Item(112, 17).Value = "12/11/2023 11:58"
$ws.Cells.Item(112, 18).Value = 3.79
$ws.Cells.Item(112, 19).Value = "09/11/2023 09:13"
$ws.Cells.Item(112, 20).Value = 5.15
$ws.Cells.Item(112, 21).Value = "12/11/2023 11:56"
$ws.Cells.Item(112, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/cordoba-ad-ceuta/KSogElxb/"

# Row 113
$ws.Cells.Item(111, 1).Copy() | Out-Null
$ws.Cells.Item(113, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111, 5).Copy() | Out-Null
$ws.Cells.Item(113, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(113, 1).Value = 112
$ws.Cells.Item(113, 2).Value = "spain"
$ws.Cells.Item(113, 3).Value = "primera-rfef-group-2"
$ws.Cells.Item(113, 4).Value = "2023-2024"
$ws.Cells.Item(113, 5).Value = 45242.5
$ws.Cells.Item(113, 6).Value = "Alcoyano"
$ws.Cells.Item(113, 7).Value = 3
$ws.Cells.Item(113, 8).Value = "Melilla"
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 1.77
$ws.Cells.Item(113, 11).Value = "09/11/2023 09:13"
$ws.Cells.Item(113, 12).Value = 1.81
$ws.Cells.Item(113, 13).Value = "12/11/2023 11:55"
$ws.Cells.Item(113, 14).Value = 3.12
$ws.Cells.Item(113, 15).Value = "09/11/2023 09:13"
$ws.Cells.Item(113, 16).Value = 3.15
$ws.Cells.Item(113, 17).Value = "12/11/2023 11:55"
$ws.Cells.Item(113, 18).Value = 4.63
$ws.Cells.Item(113, 19).Value = "09/11/2023 09:13"
$ws.Cells.Item(113, 20).Value = 5.34
$ws.Cells.Item(113, 21).Value = "12/11/2023 11:55"
$ws.Cells.Item(113, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/alcoyano-melilla/bc9M0V7b/"

# Row 114
$ws.Cells.Item(111, 1).Copy() | Out-Null
$ws.Cells.Item(114, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111, 5).Copy() | Out-Null
$ws.Cells.Item(114, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(114, 1).Value = 113
$ws.Cells.Item(114, 2).Value = "spain"
$ws.Cells.Item(114, 3).Value = "primera-rfef-group-2"
$ws.Cells.Item(114, 4).Value = "2023-2024"
$ws.Cells.Item(114, 5).Value = 45242.5
$ws.Cells.Item(114, 6).Value = "Baleares"
$ws.Cells.Item(114, 7).Value = 1
$ws.Cells.Item(114, 8).Value = "Merida AD"
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 2.25
$ws.Cells.Item(114, 11).Value = "09/11/2023 09:13"
$ws.Cells.Item(114, 12).Value = 2.46
$ws.Cells.Item(114, 13).Value = "12/11/2023 11:30"
$ws.Cells.Item(114, 14).Value = 2.94
$ws.Cells.Item(114, 15).Value = "09/11/2023 09:13"
$ws.Cells.Item(114, 16).Value = 2.94
$ws.Cells.Item(114, 17).Value = "12/11/2023 10:01"
$ws.Cells.Item(114, 18).Value = 3.26
$ws.Cells.Item(114, 19).Value = "09/11/2023 09:13"
$ws.Cells.Item(114, 20).Value = 3.23
$ws.Cells.Item(114, 21).Value = "12/11/2023 11:30"
$ws.Cells.Item(114, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/baleares-merida-ad/Sx1wcm7N/"

# Row 115
$ws.Cells.Item(111, 1).Copy() | Out-Null
$ws.Cells.Item(115, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(111, 5).Copy() | Out-Null
$ws.Cells.Item(115, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(115, 1).Value = 114
$ws.Cells.Item(115, 2).Value = "spain"
$ws.Cells.Item(115, 3).Value = "primera-rfef-group-2"
$ws.Cells.Item(115, 4).Value = "2023-2024"
$ws.Cells.Item(115, 5).Value = 45242.5
$ws.Cells.Item(115, 6).Value = "Sanluqueno"
$ws.Cells.Item(115, 7).Value = 3
$ws.Cells.Item(115, 8).Value = "Linares"
$ws.Cells.Item(115, 9).Value = 2
$ws.Cells.Item(115, 10).Value = 1.86
$ws.Cells.Item(115, 11).Value = "11/11/2023 00:12"
$ws.Cells.Item(115, 12).Value = 1.81
$ws.Cells.Item(115, 13).Value = "12/11/2023 11:45"
$ws.Cells.Item(115, 14).Value = 3.21
$ws.Cells.Item(115, 15).Value = "11/11/2023 00:12"
$ws.Cells.Item(115, 16).Value = 3.47
$ws.Cells.Item(115, 17).Value = "12/11/2023 11:45"
$ws.Cells.Item(115, 18).Value = 3.95
$ws.Cells.Item(115, 19).Value = "11/11/2023 00:12"
$ws.Cells.Item(115, 20).Value = 4.63
$ws.Cells.Item(115, 21).Value = "12/11/2023 11:45"
$ws.Cells.Item(115, 22).Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/sanluqueno-linares/x4dlFUMi/"

